# Updates the cryptos price/volume/hour table to reflect the latest GitHub Actions scrape.
# For each touched cell we force a Text number format before writing the value so that
# numeric-looking strings (prices, percentages, hour) stay stored as text, matching the
# original inline-string cell type, then reset the style back to Normal so no stray
# cell formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($Worksheet, $CellRef, $NewValue)
    $range = $Worksheet.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $NewValue
    $range.Style = "Normal"
}

$updates = @(
    @{ Cell = "D2"; Value = "306.63" }
    @{ Cell = "E2"; Value = "1.01%" }
    @{ Cell = "G2"; Value = "21" }
    @{ Cell = "D3"; Value = "35.89" }
    @{ Cell = "E3"; Value = "-0.49%" }
    @{ Cell = "G3"; Value = "21" }
    @{ Cell = "D4"; Value = "4.993" }
    @{ Cell = "E4"; Value = "-0.95%" }
    @{ Cell = "G4"; Value = "21" }
    @{ Cell = "D5"; Value = "0.08090" }
    @{ Cell = "E5"; Value = "0.10%" }
    @{ Cell = "G5"; Value = "21" }
    @{ Cell = "D6"; Value = "1.898" }
    @{ Cell = "E6"; Value = "-3.82%" }
    @{ Cell = "G6"; Value = "21" }
    @{ Cell = "D7"; Value = "4.151" }
    @{ Cell = "E7"; Value = "2.13%" }
    @{ Cell = "G7"; Value = "21" }
    @{ Cell = "D8"; Value = "7.886" }
    @{ Cell = "E8"; Value = "1.06%" }
    @{ Cell = "G8"; Value = "21" }
    @{ Cell = "D9"; Value = "0.9317" }
    @{ Cell = "E9"; Value = "0.23%" }
    @{ Cell = "G9"; Value = "21" }
    @{ Cell = "D10"; Value = "0.1243" }
    @{ Cell = "E10"; Value = "-14.98%" }
    @{ Cell = "G10"; Value = "21" }
    @{ Cell = "D11"; Value = "0.1910" }
    @{ Cell = "E11"; Value = "0.69%" }
    @{ Cell = "G11"; Value = "21" }
    @{ Cell = "D12"; Value = "0.09245" }
    @{ Cell = "E12"; Value = "3.23%" }
    @{ Cell = "G12"; Value = "21" }
    @{ Cell = "D13"; Value = "0.03511" }
    @{ Cell = "E13"; Value = "1.79%" }
    @{ Cell = "G13"; Value = "21" }
    @{ Cell = "E14"; Value = "0.62%" }
    @{ Cell = "G14"; Value = "21" }
    @{ Cell = "D15"; Value = "0.001424" }
    @{ Cell = "E15"; Value = "2.22%" }
    @{ Cell = "G15"; Value = "21" }
    @{ Cell = "D16"; Value = "0.006283" }
    @{ Cell = "E16"; Value = "8.14%" }
    @{ Cell = "G16"; Value = "21" }
    @{ Cell = "D17"; Value = "3.616" }
    @{ Cell = "E17"; Value = "2.41%" }
    @{ Cell = "G17"; Value = "21" }
    @{ Cell = "D18"; Value = "3.107" }
    @{ Cell = "E18"; Value = "9.64%" }
    @{ Cell = "G18"; Value = "21" }
    @{ Cell = "D19"; Value = "0.3444" }
    @{ Cell = "E19"; Value = "0.01%" }
    @{ Cell = "G19"; Value = "21" }
    @{ Cell = "E20"; Value = "2.36%" }
    @{ Cell = "G20"; Value = "21" }
    @{ Cell = "D21"; Value = "5.174" }
    @{ Cell = "E21"; Value = "3.03%" }
    @{ Cell = "G21"; Value = "21" }
    @{ Cell = "E22"; Value = "5.91%" }
    @{ Cell = "G22"; Value = "21" }
    @{ Cell = "D23"; Value = "0.04416" }
    @{ Cell = "E23"; Value = "-1.19%" }
    @{ Cell = "G23"; Value = "21" }
    @{ Cell = "D24"; Value = "0.001236" }
    @{ Cell = "E24"; Value = "2.50%" }
    @{ Cell = "G24"; Value = "21" }
    @{ Cell = "D25"; Value = "0.004725" }
    @{ Cell = "E25"; Value = "-1.88%" }
    @{ Cell = "G25"; Value = "21" }
    @{ Cell = "D26"; Value = "0.0001301" }
    @{ Cell = "E26"; Value = "6.23%" }
    @{ Cell = "G26"; Value = "21" }
    @{ Cell = "G27"; Value = "21" }
    @{ Cell = "G28"; Value = "21" }
    @{ Cell = "G29"; Value = "21" }
    @{ Cell = "G30"; Value = "21" }
    @{ Cell = "G31"; Value = "21" }
    @{ Cell = "G32"; Value = "21" }
    @{ Cell = "G33"; Value = "21" }
    @{ Cell = "G34"; Value = "21" }
    @{ Cell = "G35"; Value = "21" }
    @{ Cell = "G36"; Value = "21" }
    @{ Cell = "G37"; Value = "21" }
    @{ Cell = "G38"; Value = "21" }
    @{ Cell = "D39"; Value = "0.01954" }
    @{ Cell = "E39"; Value = "2.92%" }
    @{ Cell = "G39"; Value = "21" }
    @{ Cell = "D40"; Value = "0.05177" }
    @{ Cell = "E40"; Value = "8.04%" }
    @{ Cell = "G40"; Value = "21" }
    @{ Cell = "D41"; Value = "0.007560" }
    @{ Cell = "E41"; Value = "3.07%" }
    @{ Cell = "G41"; Value = "21" }
    @{ Cell = "D42"; Value = "0.01017" }
    @{ Cell = "E42"; Value = "-3.81%" }
    @{ Cell = "G42"; Value = "21" }
    @{ Cell = "D43"; Value = "0.1375" }
    @{ Cell = "E43"; Value = "1.96%" }
    @{ Cell = "G43"; Value = "21" }
    @{ Cell = "D44"; Value = "0.002102" }
    @{ Cell = "E44"; Value = "0.04%" }
    @{ Cell = "G44"; Value = "21" }
    @{ Cell = "G45"; Value = "21" }
    @{ Cell = "D46"; Value = "0.00006420" }
    @{ Cell = "E46"; Value = "3.72%" }
    @{ Cell = "G46"; Value = "21" }
    @{ Cell = "E47"; Value = "0.48%" }
    @{ Cell = "G47"; Value = "21" }
    @{ Cell = "G48"; Value = "21" }
    @{ Cell = "G49"; Value = "21" }
    @{ Cell = "E50"; Value = "0.48%" }
    @{ Cell = "G50"; Value = "21" }
    @{ Cell = "E51"; Value = "0.48%" }
    @{ Cell = "G51"; Value = "21" }
)

foreach ($update in $updates) {
    Set-TextCellValue $ws $update.Cell $update.Value
}

